# Regenerate merged AHB files
# - Rename the _old/_new diff-column headers to _FV2210/_FV2304
# - Wrap the data range in an Excel Table (Table1) with an AutoFilter
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row (row 1) from the *_old / *_new suffix scheme to the
#    *_FV2210 / *_FV2304 scheme. Column K ("diff") is unchanged.
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn the used range A1:U57 into a native Excel table ("Table1") with
#    an AutoFilter on the header row, matching the data already on the sheet.
$dataRange = $ws.Range("A1:U57")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# 3) Freeze panes so the header row (row 1) stays visible while scrolling.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit complete"
